# Applies the OFB-IdF src_idf "Reseau Loup/Lynx" sheet content refresh
# (commit 53c6ab1) to the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title cells (row 1) -----------------------------------------------
$ws.Range("C1").Value = "Réseau Loup/Lynx (loup)"
$ws.Range("L1").Value = "Réseau Loup/Lynx (loup)"

# --- Description text (row 8) -------------------------------------------
$ws.Range("C8").Value = "Espèce protégée (convention de Berne), recolonisant progressivement la France, le loup gris (Canis lupus) est encore très peu observé en Ile-de-France. Comme pour tous les territoires en front de colonisation, le réseau Loup est déployé au niveau « Sentinelle » depuis 2017, et différentes procédures peuvent alors être mises en place en cas de signalement ou de détection d’un grand canidé."

# --- Formations list re-ordered (row 9) ----------------------------------
$ws.Range("O9").Value = "Formation correspondant de réseau ( 3 jours)`nFormation constat de dommage (1,5 jours)"

# --- Roles DN / DR / SD (row 11) -----------------------------------------
$ws.Range("K11").Value = "Rapportage PNA`nAnimation & Communication`nMéthodologie`nProduction bilans`nGestion des données"
$ws.Range("L11").Value = "Centralisation et analyse des signalements`nSaisie des données`nFormation"
$ws.Range("M11").Value = "Recueil d'indices`nConstat de dommages`nCellule de veille"

# --- Objectif (row 12) ----------------------------------------------------
$ws.Range("C12").Value = "Veille`nObservations opportunistes et recueil d'indices`nConstats de dommages"

# --- Animation contacts: add the PPC correspondent (row 13) ---------------
$ws.Range("G13").Value = "Animation nationale:`nNicolas JEAN`n`nAnimation régionale:`nSamuel DEMBSKI`n`nCorrespondants départementaux:`nPPC: Arnaud LOIZE`n77: Corinne REVEL`n       Julien CURE`n78-95: Estelle DEBOST`n91: Philippe TURQUIN`n`nCourriel du réseau:`nreseau.loup-lynx@ofb.gouv.fr"

# --- Diffusion / Saisie des donnees (row 16) -------------------------------
$ws.Range("C16").Value = "Les données contribuent à la mise en oeuvre du PNA Loup notamment pour l'estimation annuelle de l'effectif de loups en France. Elles sont également valorisées dans les flash infos loup, les bilans saisonniers ou annuels et contribuent à de nombreux travaux scientifiques. Tout cela permet de mieux connaître la population de loups (aire de répartition, démographie) et de suivre son évolution afin d’accompagner les acteurs et d'aider à la mise en place de mesures de protection."
$ws.Range("K16").Value = "9 « fiches indices » liées chacune à un type d’évènement (observation visuelle, photo, empreintes/piste, excrément/poils, hurlement, cadavre de proie sauvage, cadavre de proie domestique, urine/sang, cadavre) sont disponibles et doivent être renseignées par le correspondant après entretien avec l'observateur.`nLa localisation précise (coordonnées géographiques) de l'observation est systématiquement relevée.`nStockage des prélèvements de matériel biologique dans un congélateur spécifique. Eviter la congélation/décongélation pour les analyses ADN. Délai de 48h max pour analyse de cadavre."
$ws.Range("O16").Value = "Pour chaque signalement, prendre le kit matériel adapté selon la situation (cf. matériel détaillé pour chaque cas dans le guide réflexe).`nDe manière systématique prévoir:`nfiche adaptée, GPS, appareil photo.`nPour le prélèvement de matériel biologique prévoir: `nsac de récupération d’indice, feutre indélébile, gants, scalpel, masque, gel hydroalcoolique."

# --- SAGIR procedure (row 29) ----------------------------------------------
$ws.Range("L29").Value = "Le correspondant départemental complète la ou les fiches indices adéquates à partir du témoignage de l'observateur et fournit une carte de localisation, avec si possible des coordonnées GPS et tout élément pertinent (photo/vidéo, échantillons).`nTransmettre les éléments à l’animateur régional du réseau pour expertise.`nClassement de l’évènement après analyse, et transmission du résultat à la DDT par le correspondant départemental."
$ws.Range("O29").Value = "https://www.loupfrance.fr/carte-des-indices-de-presence-transmis-au-reseau-loup-lynx/"

# --- Footer links (rows 47-49) ----------------------------------------------
# C47 now points at the general info site instead of the "Plan loup" page.
$ws.Range("C47").Formula = '=HYPERLINK("https://www.loupfrance.fr", "Site d''information")'
# L47 bumps the guide version referenced on the DR file server.
$ws.Range("L47").Formula = '=HYPERLINK("\\ad.intra\dfs\COMMUNS\REGIONS\IDF\DR\05_CONNAISSANCE\Loup\Guide réflexe réseau Loup Lynx_DRIDF_v2.4.pdf", "Guide réflexe (serveur DR)")'

# C48 used to be a blank, unformatted cell (style shared with D48:H48).
# Give it the same look as the other footer-link cells by copying L47's
# formatting, then move the former C47 "Plan loup" hyperlink into it.
$ws.Range("L47").Copy() | Out-Null
$ws.Range("C48").PasteSpecial(-4122) | Out-Null
$ws.Range("C48").Formula = '=HYPERLINK("https://agriculture.gouv.fr/plan-loup-un-nouveau-cadre-national-dactions-pour-renforcer-la-coexistence-du-loup-et-des-activites", "Plan loup")'

# --- Edited-on date + signalement-sheet link (row 49) -----------------------
$ws.Range("A49").Value = "Editée le 2025-03-14"
$ws.Range("L49").Formula = '=HYPERLINK("\\ad.intra\dfs\COMMUNS\REGIONS\IDF\DR\05_CONNAISSANCE\Loup", "Fiches indice (serveur DR)")'
